# Update workbook "horarios-141" with refreshed scrape data (10:55:35).
# Applies new Hora_Scrap / Linea / Minutos values to existing rows and
# appends newly scraped rows to each of the three sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws1.Cells.Item(2,1).Value = "Última actualización: 10:55:35"
$ws1.Cells.Item(3,1).Value = "Total filas: 155"
$ws1.Cells.Item(28,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(29,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(59,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(60,3).Value = "215A_EL PATO"
$ws1.Cells.Item(62,1).Value = "08:00:50"
$ws1.Cells.Item(62,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(62,4).Value = 3
$ws1.Cells.Item(63,1).Value = "06:46:06"
$ws1.Cells.Item(63,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(63,4).Value = 77
$ws1.Cells.Item(71,1).Value = "08:00:50"
$ws1.Cells.Item(71,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(71,4).Value = 33
$ws1.Cells.Item(72,1).Value = "08:30:59"
$ws1.Cells.Item(72,3).Value = "215C_EL PATO"
$ws1.Cells.Item(72,4).Value = 3
$ws1.Cells.Item(91,1).Value = "08:48:29"
$ws1.Cells.Item(91,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(91,4).Value = 13
$ws1.Cells.Item(92,1).Value = "08:56:14"
$ws1.Cells.Item(92,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(92,4).Value = 5
$ws1.Cells.Item(94,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(95,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(106,1).Value = "08:56:14"
$ws1.Cells.Item(106,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(106,4).Value = 37
$ws1.Cells.Item(107,1).Value = "08:30:59"
$ws1.Cells.Item(107,3).Value = "15_ABASTO"
$ws1.Cells.Item(107,4).Value = 63
$ws1.Cells.Item(126,1).Value = "10:26:25"
$ws1.Cells.Item(126,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(126,4).Value = 7
$ws1.Cells.Item(127,1).Value = "08:56:14"
$ws1.Cells.Item(127,3).Value = "14_ABASTO"
$ws1.Cells.Item(127,4).Value = 97
$ws1.Cells.Item(128,1).Value = "10:26:25"
$ws1.Cells.Item(128,3).Value = "15_ABASTO"
$ws1.Cells.Item(128,4).Value = 8
$ws1.Cells.Item(129,1).Value = "09:31:15"
$ws1.Cells.Item(129,3).Value = "14_ABASTO"
$ws1.Cells.Item(129,4).Value = 63
$ws1.Cells.Item(135,1).Value = "10:55:35"
$ws1.Cells.Item(135,4).Value = 1
$ws1.Cells.Item(136,1).Value = "10:55:35"
$ws1.Cells.Item(136,2).Value = "10:56"
$ws1.Cells.Item(136,4).Value = 1
$ws1.Cells.Item(137,2).Value = "10:57"
$ws1.Cells.Item(137,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(137,4).Value = 31
$ws1.Cells.Item(138,1).Value = "10:55:35"
$ws1.Cells.Item(138,2).Value = "11:01"
$ws1.Cells.Item(138,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(138,4).Value = 6
$ws1.Cells.Item(139,1).Value = "10:55:35"
$ws1.Cells.Item(139,2).Value = "11:03"
$ws1.Cells.Item(139,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(139,4).Value = 8
$ws1.Cells.Item(140,1).Value = "10:55:35"
$ws1.Cells.Item(140,2).Value = "11:04"
$ws1.Cells.Item(140,3).Value = "17_ROMERO"
$ws1.Cells.Item(140,4).Value = 9
$ws1.Cells.Item(141,1).Value = "10:55:35"
$ws1.Cells.Item(141,2).Value = "11:08"
$ws1.Cells.Item(141,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(141,4).Value = 13
$ws1.Cells.Item(142,1).Value = "10:55:35"
$ws1.Cells.Item(142,2).Value = "11:19"
$ws1.Cells.Item(142,3).Value = "215C_EL PATO"
$ws1.Cells.Item(142,4).Value = 24
$ws1.Cells.Item(143,1).Value = "10:55:35"
$ws1.Cells.Item(143,2).Value = "11:20"
$ws1.Cells.Item(143,4).Value = 25
$ws1.Cells.Item(144,1).Value = "09:31:15"
$ws1.Cells.Item(144,2).Value = "11:21"
$ws1.Cells.Item(144,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(144,4).Value = 110
$ws1.Cells.Item(145,1).Value = "10:55:35"
$ws1.Cells.Item(145,2).Value = "11:33"
$ws1.Cells.Item(145,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(145,4).Value = 38
$ws1.Cells.Item(146,1).Value = "10:55:35"
$ws1.Cells.Item(146,2).Value = "11:33"
$ws1.Cells.Item(146,3).Value = "215A_EL PATO"
$ws1.Cells.Item(146,4).Value = 38
$ws1.Cells.Item(147,1).Value = "10:55:35"
$ws1.Cells.Item(147,2).Value = "11:41"
$ws1.Cells.Item(147,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(147,4).Value = 46
$ws1.Cells.Item(148,1).Value = "10:55:35"
$ws1.Cells.Item(148,2).Value = "11:44"
$ws1.Cells.Item(148,3).Value = "215B_EL PATO"
$ws1.Cells.Item(148,4).Value = 49
$ws1.Cells.Item(149,1).Value = "10:55:35"
$ws1.Cells.Item(149,2).Value = "11:49"
$ws1.Cells.Item(149,3).Value = "15_ABASTO"
$ws1.Cells.Item(149,4).Value = 54
$ws1.Cells.Item(150,1).Value = "10:55:35"
$ws1.Cells.Item(150,2).Value = "11:51"
$ws1.Cells.Item(150,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(150,4).Value = 56
$ws1.Cells.Item(150,5).Value = "LP1912"
$ws1.Cells.Item(151,1).Value = "10:55:35"
$ws1.Cells.Item(151,2).Value = "11:56"
$ws1.Cells.Item(151,3).Value = "225_GOMEZ"
$ws1.Cells.Item(151,4).Value = 61
$ws1.Cells.Item(151,5).Value = "LP1912"
$ws1.Cells.Item(152,1).Value = "10:55:35"
$ws1.Cells.Item(152,2).Value = "12:04"
$ws1.Cells.Item(152,3).Value = "17_ROMERO"
$ws1.Cells.Item(152,4).Value = 69
$ws1.Cells.Item(152,5).Value = "LP1912"
$ws1.Cells.Item(153,1).Value = "10:55:35"
$ws1.Cells.Item(153,2).Value = "12:08"
$ws1.Cells.Item(153,3).Value = "14_ABASTO"
$ws1.Cells.Item(153,4).Value = 73
$ws1.Cells.Item(153,5).Value = "LP1912"
$ws1.Cells.Item(154,1).Value = "10:55:35"
$ws1.Cells.Item(154,2).Value = "12:19"
$ws1.Cells.Item(154,3).Value = "15_ABASTO"
$ws1.Cells.Item(154,4).Value = 84
$ws1.Cells.Item(154,5).Value = "LP1912"
$ws1.Cells.Item(155,1).Value = "10:55:35"
$ws1.Cells.Item(155,2).Value = "12:20"
$ws1.Cells.Item(155,3).Value = "10_OLMOS"
$ws1.Cells.Item(155,4).Value = 85
$ws1.Cells.Item(155,5).Value = "LP1912"
$ws1.Cells.Item(156,1).Value = "10:55:35"
$ws1.Cells.Item(156,2).Value = "12:32"
$ws1.Cells.Item(156,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(156,4).Value = 97
$ws1.Cells.Item(156,5).Value = "LP1912"
$ws1.Cells.Item(157,1).Value = "10:55:35"
$ws1.Cells.Item(157,2).Value = "12:34"
$ws1.Cells.Item(157,3).Value = "215C_EL PATO"
$ws1.Cells.Item(157,4).Value = 99
$ws1.Cells.Item(157,5).Value = "LP1912"
$ws1.Cells.Item(158,1).Value = "10:55:35"
$ws1.Cells.Item(158,2).Value = "12:36"
$ws1.Cells.Item(158,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(158,4).Value = 101
$ws1.Cells.Item(158,5).Value = "LP1912"
$ws1.Cells.Item(159,1).Value = "10:55:35"
$ws1.Cells.Item(159,2).Value = "12:47"
$ws1.Cells.Item(159,3).Value = "10_OLMOS"
$ws1.Cells.Item(159,4).Value = 112
$ws1.Cells.Item(159,5).Value = "LP1912"
$ws1.Cells.Item(160,1).Value = "10:55:35"
$ws1.Cells.Item(160,2).Value = "12:51"
$ws1.Cells.Item(160,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(160,4).Value = 116
$ws1.Cells.Item(160,5).Value = "LP1912"

$ws2.Cells.Item(2,1).Value = "Última actualización: 10:55:35"
$ws2.Cells.Item(3,1).Value = "Total filas: 29"
$ws2.Cells.Item(23,3).Value = "215A_EL PATO"
$ws2.Cells.Item(24,3).Value = "215B_EL PATO"
$ws2.Cells.Item(31,1).Value = "10:55:35"
$ws2.Cells.Item(31,4).Value = 24
$ws2.Cells.Item(32,1).Value = "10:55:35"
$ws2.Cells.Item(32,4).Value = 38
$ws2.Cells.Item(33,1).Value = "10:55:35"
$ws2.Cells.Item(33,4).Value = 49
$ws2.Cells.Item(34,1).Value = "10:55:35"
$ws2.Cells.Item(34,2).Value = "12:34"
$ws2.Cells.Item(34,3).Value = "215C_EL PATO"
$ws2.Cells.Item(34,4).Value = 99
$ws2.Cells.Item(34,5).Value = "LP1912"

$ws3.Cells.Item(2,1).Value = "Última actualización: 10:55:35"
$ws3.Cells.Item(3,1).Value = "Total filas: 26"
$ws3.Cells.Item(30,1).Value = "10:55:35"
$ws3.Cells.Item(30,2).Value = "11:55"
$ws3.Cells.Item(30,4).Value = 60
$ws3.Cells.Item(31,1).Value = "10:26:25"
$ws3.Cells.Item(31,2).Value = "11:56"
$ws3.Cells.Item(31,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(31,4).Value = 90
$ws3.Cells.Item(31,5).Value = "L6203"

Write-Host "Horarios actualizados Linea 141 - 1234"
